# "correction - Mathis fait"
# Fill in the grading grid: copy the "pointage" (max points) column into the
# "note" (awarded points) column for every criterion, adjusting a couple of
# rows downward, and leave a comment in column D wherever points were lost
# (plus one extra comment on a row that kept full marks).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> (note awardee value, optional comment)
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "Attention: logo de page de produit est déformé"

$ws.Range("C3").Value = 8

$ws.Range("C4").Value = 2

$ws.Range("C5").Value = 10
$ws.Range("D5").Value = "préférable de centrer les produits, au lieu de laisser 3/4 de la page vide"

$ws.Range("C6").Value = 10

$ws.Range("C7").Value = 10

$ws.Range("C8").Value = 8
$ws.Range("D8").Value = "manque une section de jobs pour George"

$ws.Range("C9").Value = 10

$ws.Range("C10").Value = 5

$ws.Range("C11").Value = 10

$ws.Range("C12").Value = 10

$ws.Range("C13").Value = 10
$ws.Range("D13").Value = "bon commentaires bien utile en html, mais manque de commentaires en css"

$ws.Range("C14").Value = 10

$ws.Range("C15").Value = 5

# Column D mirrors column A's wrap-text formatting for every data row, even
# the rows that received no comment.
$ws.Range("D2:D15").WrapText = $true

# The corrector's cursor ended up on D16 (just past the last comment cell)
# with the view scrolled back up to the top of the sheet.
$ws.Range("D16").Select() | Out-Null

Write-Output "applied correction grid"
